$d = $word.ActiveDocument

# Locate the three Heading1 "Test N" paragraphs touched by this edit by
# scanning all paragraphs for the Heading1 style and matching text, so the
# script is resilient to any paragraph-index drift.
$targets = @()
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text
    if ($txt -eq "Test 7`r" -or $txt -eq "Test 8`r") {
        $styleName = $para.Range.ParagraphStyle.NameLocal
        if ($styleName -eq "Heading 1") {
            $targets += $i
        }
    }
}

# $targets now holds, in document order: [first "Test 7", second "Test 7", "Test 8"]
$idxFirstTest7  = $targets[0]
$idxSecondTest7 = $targets[1]
$idxTest8       = $targets[2]

# --- 1) First "Test 7" heading: keep the text, but split the single run
#        "Test 7" into two runs "Test " + "7" with identical formatting
#        (mirrors the author's original run split). ---
$p1 = $d.Paragraphs.Item($idxFirstTest7)
$s1 = $p1.Range.Start
$digit1 = $d.Range($s1 + 5, $s1 + 6)
$digit1.Bold = 1
$digit1.Bold = 0

# --- 2) Second "Test 7" heading becomes "Test 8", split into "Test " + "8". ---
$p2 = $d.Paragraphs.Item($idxSecondTest7)
$s2 = $p2.Range.Start
$digit2 = $d.Range($s2 + 5, $s2 + 6)
$digit2.Text = "8"
$digit2b = $d.Range($s2 + 5, $s2 + 6)
$digit2b.Bold = 1
$digit2b.Bold = 0

# --- 3) Existing "Test 8" heading (already split "Test "/"8") becomes
#        "Test 9" - only the digit run's text changes. ---
$p3 = $d.Paragraphs.Item($idxTest8)
$s3 = $p3.Range.Start
$digit3 = $d.Range($s3 + 5, $s3 + 6)
$digit3.Text = "9"

Write-Host "Test7a=[$($p1.Range.Text)] Test7b->8=[$($p2.Range.Text)] Test8->9=[$($p3.Range.Text)]"
